$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("C4").Select() | Out-Null

$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Sheet with Errors"

$newSheet.Range("B1").Value = "field2"
$newSheet.Range("A1").Value = "field1"
$newSheet.Range("C1").Value = "result"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("C2").Formula = "=A2/B2"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = 3

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = 4

$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = 0

$newSheet.Range("A6").Value = 5
$newSheet.Range("B6").Value = 6

$newSheet.Range("C3:C6").Formula = "=A3/B3"

$newSheet.Activate()
$newSheet.Range("E5").Select() | Out-Null
